$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# --- Row 8: period headers (rolled forward by one quarter) ---
$ws.Range("D8").Value = "3 ماهه منتهی به 1399/09"
$ws.Range("E8").Value = "6 ماهه منتهی به 1399/12"
$ws.Range("F8").Value = "9 ماهه منتهی به 1400/03"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/06"
$ws.Range("H8").Value = "3 ماهه منتهی به 1400/09"
$ws.Range("I8").Value = "6 ماهه منتهی به 1400/12"
$ws.Range("J8").Value = "9 ماهه منتهی به 1401/03"
$ws.Range("K8").Value = "12 ماهه منتهی به 1401/06"
$ws.Range("L8").Value = "3 ماهه منتهی به 1401/09"
$ws.Range("M8").Value = "6 ماهه منتهی به 1401/12"

# --- Row 9: publish-date headers (rolled forward by one quarter) ---
$ws.Range("D9").Value = "1400-10-30 (2)"
$ws.Range("E9").Value = "1401-02-21 (5)"
$ws.Range("F9").Value = "1401-04-29 (3)"
$ws.Range("G9").Value = "1401-08-15 (8)"
$ws.Range("H9").Value = "1401-10-29 (2)"
$ws.Range("I9").Value = "1402-01-30 (3)"
$ws.Range("J9").Value = "1401-04-29"
$ws.Range("K9").Value = "1402-01-30 (4)"
$ws.Range("L9").Value = "1401-10-29"
$ws.Range("M9").Value = "1402-01-30"

# --- Data rows 11-27: values rolled forward one quarter, column I recomputed
#     where the pricing-algorithm change affects it, column M is the new quarter ---
# Row 11
$ws.Range("D11").Value = 1219476
$ws.Range("E11").Value = 2370125
$ws.Range("F11").Value = 4209298
$ws.Range("G11").Value = 5757814
$ws.Range("H11").Value = 1762820
$ws.Range("I11").Value = 3232505
$ws.Range("J11").Value = 5238548
$ws.Range("K11").Value = 7182656
$ws.Range("L11").Value = 2297475
$ws.Range("M11").Value = 4359899

# Row 12
$ws.Range("D12").Value = -700724
$ws.Range("E12").Value = -1409077
$ws.Range("F12").Value = -2310746
$ws.Range("G12").Value = -3314726
$ws.Range("H12").Value = -852844
$ws.Range("I12").Value = -1989060
$ws.Range("J12").Value = -3169517
$ws.Range("K12").Value = -4518294
$ws.Range("L12").Value = -1496391
$ws.Range("M12").Value = -3007934

# Row 13
$ws.Range("D13").Value = 518752
$ws.Range("E13").Value = 961048
$ws.Range("F13").Value = 1898552
$ws.Range("G13").Value = 2443088
$ws.Range("H13").Value = 909976
$ws.Range("I13").Value = 1243445
$ws.Range("J13").Value = 2069031
$ws.Range("K13").Value = 2664362
$ws.Range("L13").Value = 801084
$ws.Range("M13").Value = 1351965

# Row 14
$ws.Range("D14").Value = -64747
$ws.Range("E14").Value = -152103
$ws.Range("F14").Value = -314340
$ws.Range("G14").Value = -614077
$ws.Range("H14").Value = -97235
$ws.Range("I14").Value = -235353
$ws.Range("J14").Value = -385170
$ws.Range("K14").Value = -565326
$ws.Range("L14").Value = -218482
$ws.Range("M14").Value = -418181

# Row 15
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 0

# Row 16
$ws.Range("D16").Value = 5282
$ws.Range("E16").Value = -4650
$ws.Range("F16").Value = -7011
$ws.Range("G16").Value = -76821
$ws.Range("H16").Value = 41126
$ws.Range("I16").Value = 50305
$ws.Range("J16").Value = 54691
$ws.Range("K16").Value = 44228
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = 1203

# Row 17
$ws.Range("D17").Value = 459287
$ws.Range("E17").Value = 804295
$ws.Range("F17").Value = 1577201
$ws.Range("G17").Value = 1752190
$ws.Range("H17").Value = 853867
$ws.Range("I17").Value = 1058397
$ws.Range("J17").Value = 1738552
$ws.Range("K17").Value = 2143264
$ws.Range("L17").Value = 582602
$ws.Range("M17").Value = 934987

# Row 18
$ws.Range("D18").Value = -39711
$ws.Range("E18").Value = -81251
$ws.Range("F18").Value = -128868
$ws.Range("G18").Value = -186707
$ws.Range("H18").Value = -62500
$ws.Range("I18").Value = -51422
$ws.Range("J18").Value = -77319
$ws.Range("K18").Value = -102282
$ws.Range("L18").Value = -29437
$ws.Range("M18").Value = -87107

# Row 19
$ws.Range("D19").Value = 971
$ws.Range("E19").Value = 42198
$ws.Range("F19").Value = 59903
$ws.Range("G19").Value = 1117927
$ws.Range("H19").Value = 16190
$ws.Range("I19").Value = 213925
$ws.Range("J19").Value = 314648
$ws.Range("K19").Value = 2090990
$ws.Range("L19").Value = 93106
$ws.Range("M19").Value = 577833

# Row 20
$ws.Range("D20").Value = 420547
$ws.Range("E20").Value = 765242
$ws.Range("F20").Value = 1508236
$ws.Range("G20").Value = 2683410
$ws.Range("H20").Value = 807557
$ws.Range("I20").Value = 1220900
$ws.Range("J20").Value = 1975881
$ws.Range("K20").Value = 4131972
$ws.Range("L20").Value = 646271
$ws.Range("M20").Value = 1425713

# Row 21
$ws.Range("D21").Value = -96706
$ws.Range("E21").Value = -97134
$ws.Range("F21").Value = -269062
$ws.Range("G21").Value = -339740
$ws.Range("H21").Value = -161511
$ws.Range("I21").Value = -85614
$ws.Range("J21").Value = -187709
$ws.Range("K21").Value = -213056
$ws.Range("L21").Value = -45239
$ws.Range("M21").Value = -92845

# Row 22
$ws.Range("D22").Value = 323841
$ws.Range("E22").Value = 668108
$ws.Range("F22").Value = 1239174
$ws.Range("G22").Value = 2343670
$ws.Range("H22").Value = 646046
$ws.Range("I22").Value = 1135286
$ws.Range("J22").Value = 1788172
$ws.Range("K22").Value = 3918916
$ws.Range("L22").Value = 601032
$ws.Range("M22").Value = 1332868

# Row 23
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 70655
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 155000
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 0

# Row 24
$ws.Range("D24").Value = 323841
$ws.Range("E24").Value = 668108
$ws.Range("F24").Value = 1239174
$ws.Range("G24").Value = 2465200
$ws.Range("H24").Value = 646046
$ws.Range("I24").Value = 1205941
$ws.Range("J24").Value = 1788172
$ws.Range("K24").Value = 4073916
$ws.Range("L24").Value = 601032
$ws.Range("M24").Value = 1332868

# Row 25
$ws.Range("D25").Value = 122
$ws.Range("E25").Value = 252
$ws.Range("F25").Value = 467
$ws.Range("G25").Value = 929
$ws.Range("H25").Value = 243
$ws.Range("I25").Value = 454
$ws.Range("J25").Value = 674
$ws.Range("K25").Value = 1535
$ws.Range("L25").Value = 227
$ws.Range("M25").Value = 502

# Row 26
$ws.Range("D26").Value = 2653560
$ws.Range("E26").Value = 2653560
$ws.Range("F26").Value = 2653560
$ws.Range("G26").Value = 2653560
$ws.Range("H26").Value = 2653560
$ws.Range("I26").Value = 2653560
$ws.Range("J26").Value = 2653560
$ws.Range("K26").Value = 2653560
$ws.Range("L26").Value = 2653560
$ws.Range("M26").Value = 2653560

# Row 27
$ws.Range("D27").Value = 122
$ws.Range("E27").Value = 252
$ws.Range("F27").Value = 467
$ws.Range("G27").Value = 929
$ws.Range("H27").Value = 243
$ws.Range("I27").Value = 454
$ws.Range("J27").Value = 674
$ws.Range("K27").Value = 1535
$ws.Range("L27").Value = 227
$ws.Range("M27").Value = 502

# --- Column widths: window rolled forward one quarter ---
$ws.Columns.Item(4).ColumnWidth = 27.166666666666668  # column D
$ws.Columns.Item(5).ColumnWidth = 27.166666666666668  # column E
$ws.Columns.Item(6).ColumnWidth = 27.166666666666668  # column F
$ws.Columns.Item(7).ColumnWidth = 28.166666666666668  # column G
$ws.Columns.Item(8).ColumnWidth = 27.166666666666668  # column H
$ws.Columns.Item(9).ColumnWidth = 27.166666666666668  # column I
$ws.Columns.Item(10).ColumnWidth = 27.166666666666668  # column J
$ws.Columns.Item(11).ColumnWidth = 28.166666666666668  # column K
$ws.Columns.Item(12).ColumnWidth = 27.166666666666668  # column L
$ws.Columns.Item(13).ColumnWidth = 27.166666666666668  # column M

# --- Row heights (minor font-metric re-layout) ---
$ws.Rows.Item(2).RowHeight = 15.6
$ws.Rows.Item(5).RowHeight = 40.8
$ws.Rows.Item(6).RowHeight = 40.8
$ws.Rows.Item(8).RowHeight = 15.6